$wb = $excel.ActiveWorkbook

# --- Overview sheet: dc6c0c59 row (row 3) flips from "Ready for handoff"
# to "Handed back: in sync with en-US" for both the zh-cn and de-de columns.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: dc6c0c59 row (row 3) handback completed.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("H3").Value = "2016-08-28 20:47:50"
$wsZhCn.Range("K3").Value = "2016-08-28 20:48:06"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet: dc6c0c59 row (row 3) handback completed.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-28 20:48:13"
$wsDeDe.Range("P3").Value = ""
